$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row (data goes from row 2 through this row).
$lastRow = $ws.Cells.Item(1, 1).End(-4121).Row

for ($r = 2; $r -le $lastRow; $r++) {
    # "Förändrad" (changed/updated) date stamp in column C: bump 45184 -> 45186.
    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2() -eq 45184) {
        $cCell.Value = 45186
    }

    # Add the "Beteckning" (column A) text as the friendly-name second
    # argument to every bare HYPERLINK(...) formula in columns S..Y.
    $aVal = $ws.Cells.Item($r, 1).Value()
    for ($c = 19; $c -le 25; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.HasFormula()) {
            $f = $cell.Formula()
            if ($f.StartsWith("=HYPERLINK(") -and -not $f.Contains(",")) {
                $newf = $f.Substring(0, $f.Length - 1) + ', "' + $aVal + '")'
                $cell.Formula = $newf
            }
        }
    }
}
